$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 29.02381252314289
$ws.Range("D2").Value = -0.1261874768571047
$ws.Range("E2").Value = 0.01592327931556235

$ws.Range("C3").Value = 28.75368709650346
$ws.Range("D3").Value = -0.5963129034965462
$ws.Range("E3").Value = 0.3555890788764812

$ws.Range("C4").Value = 28.83047841649158
$ws.Range("D4").Value = -0.5395215835084208
$ws.Range("E4").Value = 0.2910835390714339

$ws.Range("C5").Value = 29.05125791563309
$ws.Range("D5").Value = -0.4887420843669119
$ws.Range("E5").Value = 0.2388688250313137

$ws.Range("C6").Value = 29.5158752386351
$ws.Range("D6").Value = -0.03412476136489673
$ws.Range("E6").Value = 0.001164499338211148

$ws.Range("C7").Value = 29.74304293528138
$ws.Range("D7").Value = -0.006957064718619677
$ws.Range("E7").Value = 0.00004840074949906268

$ws.Range("C8").Value = 29.58077201361602
$ws.Range("D8").Value = -0.2592279863839799
$ws.Range("E8").Value = 0.06719914892469288

$ws.Range("C9").Value = 29.68211377194293
$ws.Range("D9").Value = -0.1278862280570685
$ws.Range("E9").Value = 0.01635488732666454

$ws.Range("C10").Value = 29.54752795107213
$ws.Range("D10").Value = -0.3724720489278681
$ws.Range("E10").Value = 0.1387354272325242

$ws.Range("C11").Value = 29.75683841497857
$ws.Range("D11").Value = -0.2231615850214261
$ws.Range("E11").Value = 0.0498010930292752

$ws.Range("C12").Value = 29.73810743862573
$ws.Range("D12").Value = -0.3018925613742667
$ws.Range("E12").Value = 0.09113911861311537

$ws.Range("C13").Value = 29.76147510426867
$ws.Range("D13").Value = -0.4485248957313317
$ws.Range("E13").Value = 0.201174582090802

$ws.Range("C14").Value = 29.9606812511153
$ws.Range("D14").Value = -0.2593187488846951
$ws.Range("E14").Value = 0.06724621352312354

$ws.Range("C15").Value = 29.87607033507794
$ws.Range("D15").Value = -0.5039296649220582
$ws.Range("E15").Value = 0.2539451071884578

$ws.Range("C16").Value = 29.93461129704251
$ws.Range("D16").Value = -0.5053887029574931
$ws.Range("E16").Value = 0.2554177410770572

$ws.Range("C17").Value = 30.32728585566579
$ws.Range("D17").Value = -0.1527141443342082
$ws.Range("E17").Value = 0.02332160987972936

$ws.Range("C18").Value = 30.28804325780069
$ws.Range("D18").Value = -0.4019567421993067
$ws.Range("E18").Value = 0.1615692225994799

$ws.Range("C19").Value = 31.35097886077029
$ws.Range("D19").Value = 0.6009788607702902
$ws.Range("E19").Value = 0.3611755910927558

$ws.Range("C20").Value = 31.49617511377952
$ws.Range("D20").Value = 0.5561751137795206
$ws.Range("E20").Value = 0.3093307571876626

$ws.Range("C21").Value = 31.79926696385362
$ws.Range("D21").Value = 0.8492669638536192
$ws.Range("E21").Value = 0.7212543758931446

$ws.Range("C22").Value = 31.36317758864411
$ws.Range("D22").Value = 0.3431775886441137
$ws.Range("E22").Value = 0.1177708573475885

$ws.Range("C23").Value = 31.1702265847368
$ws.Range("D23").Value = 0.05022658473680153
$ws.Range("E23").Value = 0.002522709814323104

$ws.Range("C24").Value = 31.46856948383135
$ws.Range("D24").Value = 0.1885694838313441
$ws.Range("E24").Value = 0.03555845023241954

$ws.Range("C25").Value = 31.44570478500628
$ws.Range("D25").Value = 0.06570478500628596
$ws.Range("E25").Value = 0.00431711877272226

$ws.Range("C26").Value = 31.15478729408748
$ws.Range("D26").Value = -0.4252127059125144
$ws.Range("E26").Value = 0.1808058452694424

$ws.Range("C27").Value = 31.58687360071157
$ws.Range("D27").Value = -0.06312639928842856
$ws.Range("E27").Value = 0.003984942287122113

$ws.Range("C28").Value = 31.54552549821508
$ws.Range("D28").Value = -0.3344745017849213
$ws.Range("E28").Value = 0.1118731923442713

$ws.Range("C29").Value = 31.93445834616257
$ws.Range("D29").Value = -0.3455416538374294
$ws.Range("E29").Value = 0.1193990345367059

$ws.Range("C30").Value = 32.3731268104054
$ws.Range("D30").Value = -0.07687318959460754
$ws.Range("E30").Value = 0.005909487278448477

$ws.Range("C31").Value = 33.32158228747896
$ws.Range("D31").Value = 0.4715822874789595
$ws.Range("E31").Value = 0.222389853863888

$ws.Range("C32").Value = 33.54253931907255
$ws.Range("D32").Value = 0.6425393190725472
$ws.Range("E32").Value = 0.4128567765542127

$ws.Range("C33").Value = 34.09321692438535
$ws.Range("D33").Value = 0.9932169243853508
$ws.Range("E33").Value = 0.9864798588854957

$ws.Range("C34").Value = 33.39833273060193
$ws.Range("D34").Value = -0.001667269398069493
$ws.Range("E34").Value = 0.00000277978724573901

$ws.Range("C35").Value = 34.39789626642995
$ws.Range("D35").Value = 0.6978962664299502
$ws.Range("E35").Value = 0.487059198696864

$ws.Range("C36").Value = 35.00093290383575
$ws.Range("D36").Value = 0.900932903835745
$ws.Range("E36").Value = 0.8116800972139077

$ws.Range("C37").Value = 35.30233051941863
$ws.Range("D37").Value = 0.9023305194186335
$ws.Range("E37").Value = 0.8142003662743009

$ws.Range("C38").Value = 34.95301991141709
$ws.Range("D38").Value = 0.05301991141708839
$ws.Range("E38").Value = 0.0028111110066759

$ws.Range("C39").Value = 34.72306620000119
$ws.Range("D39").Value = -0.5769337999988053
$ws.Range("E39").Value = 0.3328526095810615

$ws.Range("C40").Value = 35.16858681465516
$ws.Range("D40").Value = -0.5314131853448387
$ws.Range("E40").Value = 0.2823999735583479

$ws.Range("C41").Value = 35.61014068077346
$ws.Range("D41").Value = -0.6898593192265352
$ws.Range("E41").Value = 0.4759058803236986

$ws.Range("C42").Value = 35.68210540268343
$ws.Range("D42").Value = -1.117894597316571
$ws.Range("E42").Value = 1.249688330709579

$ws.Range("C43").Value = 36.89876293834338
$ws.Range("D43").Value = -0.4012370616566159
$ws.Range("E43").Value = 0.160991179646835

$ws.Range("C44").Value = 37.87543651653033
$ws.Range("D44").Value = -0.02456348346967019
$ws.Range("E44").Value = 0.0006033647201647606

$ws.Range("C45").Value = 38.58766495887683
$ws.Range("D45").Value = 0.08766495887682879
$ws.Range("E45").Value = 0.007685145014876083

$ws.Range("C46").Value = 39.10540492532465
$ws.Range("D46").Value = 0.2054049253246504
$ws.Range("E46").Value = 0.0421911833476252

$ws.Range("C47").Value = 40.05100399181302
$ws.Range("D47").Value = 0.6510039918130204
$ws.Range("E47").Value = 0.4238061973564871

$ws.Range("C48").Value = 40.56251715606416
$ws.Range("D48").Value = 0.6625171560641618
$ws.Range("E48").Value = 0.438928982079345

$ws.Range("C49").Value = 38.33997861283432
$ws.Range("D49").Value = -1.760021387165686
$ws.Range("E49").Value = 3.097675283280626

$ws.Range("C50").Value = 40.43296999164407
$ws.Range("D50").Value = -0.1670300083559297
$ws.Range("E50").Value = 0.02789902369138196

$ws.Range("C51").Value = 40.02033061713695
$ws.Range("D51").Value = -0.8796693828630495
$ws.Range("E51").Value = 0.7738182231466584

$ws.Range("C52").Value = -3.821628583580964
$ws.Range("E52").Value = 15.25440955466331

$ws.Range("E53").Value = 0.3050881910932662
